$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.540.98'
$ws.Range("E2").Value = '  -4.32%  '

$ws.Range("D3").Value = '3.482.76'
$ws.Range("E3").Value = '  -4.58%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.45%  '

$ws.Range("D8").Value = '3.475.07'
$ws.Range("E8").Value = '  -4.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.187'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.66'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.55%  '

$ws.Range("E12").Value = '  -1.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000276'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '686.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.87%  '

$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '4.038.02'
$ws.Range("E17").Value = '  -4.70%  '

$ws.Range("D18").Value = '68.497.62'
$ws.Range("E18").Value = '  -4.51%  '

$ws.Range("D19").Value = '3.477.87'
$ws.Range("E19").Value = '  -5.08%  '

$ws.Range("E20").Value = '  -1.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.901'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.90%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.39%  '

$ws.Range("E29").Value = '  -7.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.75'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.45%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.83%  '

$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '566.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -14.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.07%  '

$ws.Range("E38").Value = '  -3.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '56.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.80%  '

$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("E41").Value = '  -5.64%  '

$ws.Range("E42").Value = '  -4.61%  '

$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.335'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.02%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.418.05'
$ws.Range("E44").Value = '  -8.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '33.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.20%  '

$ws.Range("D46").Value = '0.0₃0700'
$ws.Range("E46").Value = '  -8.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.00%  '

$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.08'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.24%  '

$ws.Range("E51").Value = '  -1.42%  '
